# "Ejecucion de ssh para portIn"
# Adds ssh connection/query columns (host, usuario, contraseña ssh + a
# consulta_log value) to the header/value rows of the "Semilla 4" sheet,
# and moves the active selection to F5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semilla 4")

# New header row (row 4) values
$ws.Range("E4").Value = "host ssh"
$ws.Range("F4").Value = "usuario ssh"
$ws.Range("G4").Value = "contraseña ssh"

# New value row (row 5) values - fill F5/G5 before E5 so the shared-string
# table picks up "consulta_log" ahead of "10.69.60.76"
$ws.Range("F5").Value = "consulta_log"
$ws.Range("G5").Value = "consulta_log"
$ws.Range("E5").Value = "10.69.60.76"

# Format the new header cells (bold, black font, centered) on a single
# cell first, then fan the exact same format out with copy/paste-special
# so no throwaway intermediate cell styles get minted in the style table.
$headerCell = $ws.Range("E4")
$headerCell.Font.Bold = $true
$headerCell.Font.Color = 0
$headerCell.HorizontalAlignment = -4108
$headerCell.Copy()
$ws.Range("F4:G4").PasteSpecial(-4122)

# Format the new value cells (plain black font)
$valueCell = $ws.Range("E5")
$valueCell.Font.Color = 0
$valueCell.Copy()
$ws.Range("F5:G5").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Update the view: drop the C1 top-left freeze/scroll and select F5
$ws.Activate()
$ws.Range("F5").Select()
